$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell with the default (unstyled) style. Used to strip the
# quote-prefix styling Excel applies when a numeric-looking string is
# forced to stay text, so the cell keeps the same style as the original
# file (all Price-column cells are stored as text, not numbers).
$plainStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '57.287.75'
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").Value = '3.011.65'
$ws.Range("E3").Value = '  +0.05%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = '''' + '511.44'
$ws.Range("D5").Style = $plainStyle
$ws.Range("E5").Value = '  +0.39%  '
$ws.Range("D6").Value = '''' + '138.87'
$ws.Range("D6").Style = $plainStyle
$ws.Range("E6").Value = '  +1.44%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '''' + '0.437'
$ws.Range("D8").Style = $plainStyle
$ws.Range("E8").Value = '  +0.77%  '
$ws.Range("D9").Value = '''' + '7.52'
$ws.Range("D9").Style = $plainStyle
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("D10").Value = '''' + '0.110'
$ws.Range("D10").Style = $plainStyle
$ws.Range("E10").Value = '  +0.86%  '
$ws.Range("D11").Value = '''' + '0.366'
$ws.Range("D11").Style = $plainStyle
$ws.Range("E11").Value = '  +3.35%  '
$ws.Range("D12").Value = '3.544.07'
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("E13").Value = '  +1.51%  '
$ws.Range("D14").Value = '''' + '26.51'
$ws.Range("D14").Style = $plainStyle
$ws.Range("E14").Value = '  +2.83%  '
$ws.Range("D15").Value = '''' + '0.0000165'
$ws.Range("D15").Style = $plainStyle
$ws.Range("E15").Value = '  +6.98%  '
$ws.Range("D16").Value = '57.437.16'
$ws.Range("E16").Value = '  +1.39%  '
$ws.Range("D17").Value = '''' + '6.22'
$ws.Range("D17").Style = $plainStyle
$ws.Range("E17").Value = '  +5.80%  '
$ws.Range("D18").Value = '3.020.06'
$ws.Range("E18").Value = '  +0.25%  '
$ws.Range("D19").Value = '''' + '12.76'
$ws.Range("D19").Style = $plainStyle
$ws.Range("E19").Value = '  +1.76%  '
$ws.Range("D20").Value = '''' + '7.96'
$ws.Range("D20").Style = $plainStyle
$ws.Range("E20").Value = '  +0.99%  '
$ws.Range("D21").Value = '''' + '330.00'
$ws.Range("D21").Style = $plainStyle
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = '''' + '0.998'
$ws.Range("D22").Style = $plainStyle
$ws.Range("E22").Value = '  -0.15%  '
$ws.Range("D23").Value = '''' + '0.497'
$ws.Range("D23").Style = $plainStyle
$ws.Range("E23").Value = '  +3.86%  '
$ws.Range("D24").Value = '''' + '64.55'
$ws.Range("D24").Style = $plainStyle
$ws.Range("E24").Value = '  +3.23%  '
$ws.Range("D25").Value = '''' + '0.168'
$ws.Range("D25").Style = $plainStyle
$ws.Range("E25").Value = '  +1.22%  '
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("D27").Value = '0.0₃0922'
$ws.Range("E27").Value = '  +0.27%  '
$ws.Range("D28").Value = '''' + '6.76'
$ws.Range("D28").Style = $plainStyle
$ws.Range("E28").Value = '  +2.85%  '
$ws.Range("D29").Value = '''' + '7.46'
$ws.Range("D29").Style = $plainStyle
$ws.Range("E29").Value = '  +6.91%  '
$ws.Range("D30").Value = '''' + '1.81'
$ws.Range("D30").Style = $plainStyle
$ws.Range("D31").Value = '''' + '1.20'
$ws.Range("D31").Style = $plainStyle
$ws.Range("E31").Value = '  -3.81%  '
$ws.Range("D32").Value = '''' + '20.61'
$ws.Range("D32").Style = $plainStyle
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("D33").Value = '''' + '4.71'
$ws.Range("D33").Style = $plainStyle
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("D34").Value = '''' + '154.71'
$ws.Range("D34").Style = $plainStyle
$ws.Range("E34").Value = '  -0.93%  '
$ws.Range("D35").Value = '''' + '5.88'
$ws.Range("D35").Style = $plainStyle
$ws.Range("E35").Value = '  +4.70%  '
$ws.Range("D36").Value = '''' + '1.27'
$ws.Range("D36").Style = $plainStyle
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").Value = '''' + '24.45'
$ws.Range("D37").Style = $plainStyle
$ws.Range("E37").Value = '  +2.77%  '
$ws.Range("D38").Value = '''' + '0.0681'
$ws.Range("D38").Style = $plainStyle
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("D39").Value = '3.053.43'
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").Value = '''' + '37.33'
$ws.Range("D40").Style = $plainStyle
$ws.Range("E40").Value = '  +2.02%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''' + '1.00'
$ws.Range("D41").Style = $plainStyle
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").Value = '''' + '3.85'
$ws.Range("D42").Style = $plainStyle
$ws.Range("E42").Value = '  +6.46%  '
$ws.Range("D43").Value = '2.300.70'
$ws.Range("E43").Value = '  +1.44%  '
$ws.Range("D44").Value = '''' + '0.652'
$ws.Range("D44").Style = $plainStyle
$ws.Range("E44").Value = '  +0.94%  '
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").Value = '''' + '0.981'
$ws.Range("D46").Style = $plainStyle
$ws.Range("E46").Value = '  -2.12%  '
$ws.Range("D47").Value = '''' + '6.03'
$ws.Range("D47").Style = $plainStyle
$ws.Range("E47").Value = '  +3.75%  '
$ws.Range("E48").Value = '  +1.05%  '
$ws.Range("D49").Value = '''' + '19.48'
$ws.Range("D49").Style = $plainStyle
$ws.Range("E49").Value = '  +1.44%  '
$ws.Range("E50").Value = '  -7.34%  '
$ws.Range("D51").Value = '''' + '0.0891'
$ws.Range("D51").Style = $plainStyle
$ws.Range("E51").Value = '  +1.70%  '
